$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the stray 0.5 value from I14 (keep existing style/border)
# ------------------------------------------------------------------
$ws.Range("I14").ClearContents()

# ------------------------------------------------------------------
# 2. "Legends:" section header on row 24 (E24:P24), styled like the
#    big bold/underlined title font used in A1:A3 (Arial 15 B/U)
# ------------------------------------------------------------------
$ws.Range("E24").Value = "Legends:"
$ws.Range("F24:P24").Value = " "
$legendHeader = $ws.Range("E24:P24")
$legendHeader.Font.Name = "Arial"
$legendHeader.Font.Size = 15
$legendHeader.Font.Bold = $true
$legendHeader.Font.Underline = $true
$ws.Range("E24:P24").Merge()

# ------------------------------------------------------------------
# 3. Legend entry 1 (rows 25-26) - blue swatch + remark/request text
# ------------------------------------------------------------------
$ws.Range("E25").Value = " "
$ws.Range("E26").Value = " "
$ws.Range("E25:E26").Interior.Color = 13411113
$ws.Range("E25:E26").Merge()

$ws.Range("F25").Value = "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."
$ws.Range("G25:P25").Value = " "
$ws.Range("F26:P26").Value = " "
$legendText1 = $ws.Range("F25:P26")
$legendText1.Font.Name = "Arial"
$legendText1.Font.Size = 11
$legendText1.Font.Bold = $true
$legendText1.Font.Underline = $true
$ws.Range("F25:P26").Merge()

# ------------------------------------------------------------------
# 4. Legend entry 2 (rows 27-28) - yellow swatch + half-day text
# ------------------------------------------------------------------
$ws.Range("E27").Value = " "
$ws.Range("E28").Value = " "
$ws.Range("E27:E28").Interior.Color = 6737151
$ws.Range("E27:E28").Merge()

$ws.Range("F27").Value = "Employee is considered half-day because of his time-in or time-out."
$ws.Range("G27:P27").Value = " "
$ws.Range("F28:P28").Value = " "
$legendText2 = $ws.Range("F27:P28")
$legendText2.Font.Name = "Arial"
$legendText2.Font.Size = 11
$legendText2.Font.Bold = $true
$legendText2.Font.Underline = $true
$ws.Range("F27:P28").Merge()

# ------------------------------------------------------------------
# 5. Legend entry 3 (rows 29-30) - red swatch + absent text
# ------------------------------------------------------------------
$ws.Range("E29").Value = " "
$ws.Range("E30").Value = " "
$ws.Range("E29:E30").Interior.Color = 6184671
$ws.Range("E29:E30").Merge()

$ws.Range("F29").Value = "Employee has no time-in and therefore, considered as absent."
$ws.Range("G29:P29").Value = " "
$ws.Range("F30:P30").Value = " "
$legendText3 = $ws.Range("F29:P30")
$legendText3.Font.Name = "Arial"
$legendText3.Font.Size = 11
$legendText3.Font.Bold = $true
$legendText3.Font.Underline = $true
$ws.Range("F29:P30").Merge()

Write-Host "Legends applied"
